$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.480.07"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "1.850.05"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'260.46"
$ws.Range("E5").Value = "  -7.85%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.5157"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "'0.3254"
$ws.Range("E8").Value = "  -7.94%  "
$ws.Range("D9").Value = "'0.06766"
$ws.Range("E9").Value = "  -5.46%  "
$ws.Range("D10").Value = "'18.89"
$ws.Range("E10").Value = "  -6.80%  "
$ws.Range("D11").Value = "'0.7720"
$ws.Range("E11").Value = "  -5.95%  "
$ws.Range("D12").Value = "1.917.16"
$ws.Range("E12").Value = "  +2.61%  "
$ws.Range("D13").Value = "'0.07716"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "'88.63"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "'5.031"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "'14.10"
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'0.000007907"
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("D20").Value = "26.503.16"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "2.089.93"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").Value = "'4.525"
$ws.Range("E22").Value = "  -5.72%  "
$ws.Range("D23").Value = "'9.530"
$ws.Range("E23").Value = "  -6.36%  "
$ws.Range("D24").Value = "'5.932"
$ws.Range("E24").Value = "  -5.12%  "
$ws.Range("D25").Value = "'2.354"
$ws.Range("E25").Value = "  -3.06%  "
$ws.Range("D26").Value = "'144.54"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").Value = "'1.655"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("E28").Value = "  -2.68%  "
$ws.Range("D29").Value = "'111.22"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "'4.205"
$ws.Range("E30").Value = "  -4.78%  "
$ws.Range("D31").Value = "'4.179"
$ws.Range("E31").Value = "  -4.36%  "
$ws.Range("D32").Value = "'0.08761"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").Value = "'0.04820"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("D34").Value = "'1.135"
$ws.Range("E34").Value = "  -3.85%  "
$ws.Range("D35").Value = "'2.841"
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").Value = "'0.6892"
$ws.Range("E36").Value = "  -7.95%  "
$ws.Range("D37").Value = "'3.116"
$ws.Range("E37").Value = "  -5.49%  "
$ws.Range("D38").Value = "'0.01804"
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").Value = "'2.212"
$ws.Range("E39").Value = "  -8.62%  "
$ws.Range("D40").Value = "'0.4908"
$ws.Range("E40").Value = "  -7.60%  "
$ws.Range("D41").Value = "'113.10"
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("D42").Value = "'0.9006"
$ws.Range("E42").Value = "  -7.81%  "
$ws.Range("D43").Value = "'6.136"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "'7.780"
$ws.Range("E45").Value = "  -5.03%  "
$ws.Range("D46").Value = "'0.4231"
$ws.Range("E46").Value = "  -8.42%  "
$ws.Range("E47").Value = "  -7.18%  "
$ws.Range("D48").Value = "'9.129"
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("D49").Value = "'0.05894"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "'35.08"
$ws.Range("E50").Value = "  -4.37%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.416"
$ws.Range("E51").Value = "  -6.71%  "
